$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '55.676.01'
$ws.Range('E2').Value = '  -2.09%  '
$ws.Range('D3').Value = '2.335.84'
$ws.Range('E3').Value = '  -2.48%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''502.05'
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('D6').Value = '''128.32'
$ws.Range('E6').Value = '  -3.31%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -2.95%  '
$ws.Range('D9').Value = '2.342.69'
$ws.Range('E9').Value = '  -2.48%  '
$ws.Range('E10').Value = '  +0.08%  '
$ws.Range('E11').Value = '  -0.17%  '
$ws.Range('D12').Value = '''4.78'
$ws.Range('E12').Value = '  +3.10%  '
$ws.Range('E13').Value = '  -1.65%  '
$ws.Range('D14').Value = '2.754.93'
$ws.Range('E14').Value = '  -2.29%  '
$ws.Range('D15').Value = '55.670.88'
$ws.Range('E15').Value = '  -1.99%  '
$ws.Range('D16').Value = '''21.51'
$ws.Range('E16').Value = '  -1.00%  '
$ws.Range('E17').Value = '  -2.03%  '
$ws.Range('D18').Value = '2.344.37'
$ws.Range('E18').Value = '  -3.63%  '
$ws.Range('E19').Value = '  -3.00%  '
$ws.Range('D20').Value = '''309.40'
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('E21').Value = '  -2.03%  '
$ws.Range('D22').Value = '''6.18'
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').Value = '''65.22'
$ws.Range('D25').Value = '''0.998'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  -1.24%  '
$ws.Range('E27').Value = '  -3.27%  '
$ws.Range('D28').Value = '''7.07'
$ws.Range('E28').Value = '  -4.45%  '
$ws.Range('D29').Value = '''171.68'
$ws.Range('E29').Value = '  -2.19%  '
$ws.Range('E30').Value = '  -1.11%  '
$ws.Range('D31').Value = '0.0₃0700'
$ws.Range('E31').Value = '  -3.13%  '
$ws.Range('E32').Value = '  -0.03%  '
$ws.Range('E33').Value = '  -1.55%  '
$ws.Range('D34').Value = '''0.998'
$ws.Range('E34').Value = '  +0.13%  '
$ws.Range('E35').Value = '  -5.48%  '
$ws.Range('E36').Value = '  -1.92%  '
$ws.Range('E37').Value = '  -2.34%  '
$ws.Range('E38').Value = '  -4.72%  '
$ws.Range('D39').Value = '''0.818'
$ws.Range('E39').Value = '  -0.79%  '
$ws.Range('D40').Value = '''36.05'
$ws.Range('E40').Value = '  -2.19%  '
$ws.Range('D41').Value = '''1.37'
$ws.Range('E41').Value = '  -4.51%  '
$ws.Range('E42').Value = '  -0.97%  '
$ws.Range('D43').Value = '''126.26'
$ws.Range('E43').Value = '  -3.81%  '
$ws.Range('D44').Value = '''4.68'
$ws.Range('E44').Value = '  -3.58%  '
$ws.Range('E45').Value = '  -2.68%  '
$ws.Range('D46').Value = '''0.0892'
$ws.Range('E46').Value = '  -2.36%  '
$ws.Range('D47').Value = '''236.17'
$ws.Range('E47').Value = '  -5.66%  '
$ws.Range('D48').Value = '''0.0474'
$ws.Range('E48').Value = '  -2.66%  '
$ws.Range('D49').Value = '''0.0204'
$ws.Range('E49').Value = '  -2.75%  '
$ws.Range('D50').Value = '''16.73'
$ws.Range('E50').Value = '  -1.79%  '
$ws.Range('E51').Value = '  +0.02%  '
